$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 26 ("كالونا") ratio column H26 changes from "0:0" to "2:0"
$ws.Range("H26").Value = "2:0"

# 2) Insert two new blank data rows at 27-28 (old row 27 total -> 29, old row 28 footer -> 30)
$ws.Range("27:28").Insert()

# 3) Build new row 27 by copying the formatting of row 26 (a fully-styled data row)
#    then overwrite with the new item's values.
$ws.Range("A26:Q26").Copy()
$ws.Range("A27:Q27").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows.Item(27).RowHeight = 25.5

$ws.Range("A27").Value = 21
$ws.Range("C27").Value = "قصافات اطفال"
$ws.Range("H27").Value = "5:0"
$ws.Range("L27").Value = "0"
$ws.Range("N27").Value = "10.00"
$ws.Range("P27").Value = "10.0000"
$ws.Range("Q27").Value = "1:0"

$ws.Range("A27:B27").Merge()
$ws.Range("C27:G27").Merge()
$ws.Range("H27:K27").Merge()
$ws.Range("L27:M27").Merge()
$ws.Range("N27:O27").Merge()

# 4) Build new row 28 the same way, for the repeated "كالونا" line.
$ws.Range("A26:Q26").Copy()
$ws.Range("A28:Q28").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows.Item(28).RowHeight = 24.75

$ws.Range("A28").Value = 22
$ws.Range("C28").Value = "كالونا "
$ws.Range("H28").Value = "0:0"
$ws.Range("L28").Value = "0"
$ws.Range("N28").Value = "15.00"
$ws.Range("P28").Value = "15.0000"
$ws.Range("Q28").Value = "1:0"

$ws.Range("A28:B28").Merge()
$ws.Range("C28:G28").Merge()
$ws.Range("H28:K28").Merge()
$ws.Range("L28:M28").Merge()
$ws.Range("N28:O28").Merge()

# 5) Update the grand-total cell (old P27, now shifted to P29) from 1194 to 1219
$ws.Range("P29").Value = 1219

# 6) Refresh the printed timestamp in the footer (now row 30) to the new save time
$ws.Range("A30").Value = "Sunday, 13 July, 2025 12:32 PM"
